# "fixed stationID not unique key bug"
#
# The Stations sheet keyed rows by StationID alone, which isn't unique
# across different uploaded files. Add a "File" column as the new first
# column on the Stations sheet (reusing the existing "File" shared string
# already used on the summary sheet) so rows can be disambiguated by
# file + StationID.

$wb = $excel.ActiveWorkbook

# --- Stations sheet: insert the new "File" column at the front ---
$stations = $wb.Worksheets.Item("Stations")

# Shift everything one column to the right and create a blank column A.
$stations.Columns.Item(1).Insert(-4161)

# Give the new A1 header the same look as the rest of the header row,
# then set its caption.
$stations.Range("B1").Copy()
$stations.Range("A1").PasteSpecial(-4122)
$stations.Range("A1").Value = "File"

# The two "O"-highlight conditional format columns shifted from Y:Z to
# Z:AA; collapse the rule back onto the (now contiguous) pair of columns.
$fc = $stations.Cells.FormatConditions.Item(1)
$fc.ModifyAppliesToRange($stations.Range("Z2:AA104857"))

# --- summary sheet: highlight the three temperature/GPS group headers ---
$summary = $wb.Worksheets.Item("summary")
$summary.Range("L1,O1,R1").Interior.ThemeColor = 8
$summary.Range("G3").Select()

# Stations becomes the active tab, with a fresh selection (selecting a
# range on a sheet activates that sheet, so this must happen last).
$stations.Activate()
$stations.Range("C5").Select()
